$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. '212.26') need to be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values (losing the original text semantics used throughout this sheet).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.965.61'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.634.32'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '212.26'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("E6").Value = '  -0.44%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '23.51'
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("E9").Value = '  -2.15%  '
$ws.Range("E10").Value = '  -0.37%  '
$ws.Range("D11").Value = '0.0881'
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '1.866.07'
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '1.637.49'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").Value = '0.562'
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").Value = '65.63'
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").Value = '27.963.45'
$ws.Range("E17").Value = '  +0.26%  '
$ws.Range("D18").Value = '232.21'
$ws.Range("E18").Value = '  +0.83%  '
$ws.Range("D19").Value = '0.0₃0725'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = '7.57'
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("D22").Value = '10.43'
$ws.Range("E22").Value = '  -4.82%  '
$ws.Range("D23").Value = '4.36'
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("E24").Value = '  -3.32%  '
$ws.Range("D25").Value = '154.75'
$ws.Range("E25").Value = '  +1.70%  '
$ws.Range("D26").Value = '6.94'
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").Value = '15.65'
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  -0.39%  '
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '3.09'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("B34").Value = 'Maker'
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D34").Value = '1.408.25'
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("E36").Value = '  +9.38%  '
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").Value = '0.870'
$ws.Range("E40").Value = '  -2.14%  '
$ws.Range("E41").Value = '  -1.16%  '
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").Value = '67.13'
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.47'
$ws.Range("E44").Value = '  +0.31%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '1.82'
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").Value = '1.776.43'
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").Value = '88.00'
$ws.Range("E48").Value = '  -1.10%  '
$ws.Range("D49").Value = '0.100'
$ws.Range("E49").Value = '  -0.75%  '
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = '7.56'
$ws.Range("E51").Value = '  -1.74%  '

# Restore the default cell style on those cells so no stray number format
# remains applied (matches original workbook which had no explicit style here).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
